$d = $word.ActiveDocument

$rng = $d.Range(0, 157)
$rng.Text = "The reason why I am limiting the race of my sample to only black and white defendants is because they account for 1.1 out of 1.2 million observations. "
